$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocabulary entries added for 2020-12-03 (14 days prior writing cycle)
$newWords = @(
  @("消息", "1.) n. information, news"),
  @("热烈", "1.) adj. warm [welcome], enthusiastic, ardent"),
  @("褒贬不一", "1.) n. phr. mixed reviews"),
  @("一举两得", "1.) idiom. two birds with one stone"),
  @("鼓励", "1.) v. encourage"),
  @("则", "1.) 就; 2.) 却"),
  @("参考", "1.) n. consultation, reference; 2.) v. consult, refer"),
  @("侵犯", "1.) v. infringe on, encroach on, violate"),
  @("涉嫌", "1.) adj. to be suspsect [in a crime], to be suspected of"),
  @("歧视", "1.) v. discriminate, look down upon; 2.) n. discrimination"),
  @("人士", "1.) n. person, figure, public figure"),
  @("亦或", "1.) conj. or"),
  @("焦虑", "1.) n. anxiety, worry, apprehension"),
  @("恶意", "1.) n. malice, spite, evil intention"),
  @("营销", "1.) n. marketing"),
  @("值得", "1.) adj. to be worth, to deserve"),
  @("提倡", "1.) n. promotion, advocacy; 2.) v. promote, advocate")
)

$newDate = "2020-12-03"
$startRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count

# Write column A (FOREIGN) for every new row first, then column B (ENGLISH),
# then column C (DATE) last -- matches the order new shared strings were
# interned in the original author's edit (all words, then all meanings,
# then the date at the very end).
for ($i = 0; $i -lt $newWords.Count; $i++) {
  $r = $startRow + $i
  $pair = $newWords[$i]
  $ws.Cells.Item($r, 1).Value = $pair[0]
}

for ($i = 0; $i -lt $newWords.Count; $i++) {
  $r = $startRow + $i
  $pair = $newWords[$i]
  $ws.Cells.Item($r, 2).Value = $pair[1]
}

for ($i = 0; $i -lt $newWords.Count; $i++) {
  $r = $startRow + $i
  $dateCell = $ws.Cells.Item($r, 3)
  $dateCell.Value = "'" + $newDate
  $dateCell.Style = "Normal"
}
